# Updated TPM-derived NATMI edge metrics for Vcan-Itgb1 (Sheet1, rows 2-17).
# Built from the authoritative cell-level diff; each row below lists only the
# columns that actually changed for that row (columns A-D, K, L are untouched).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row=2; Cells=@{ "E"="3"; "F"="1"; "G"="4.204118999999999"; "H"="12.612357"; "I"="0.01983154129720676"; "J"="0.01983154129720676"; "M"="168.1098273333333"; "N"="504.329482"; "O"="0.2984182258032519"; "P"="0.298418225803252"; "Q"="706.7537191787859"; "R"="6360.783472609073"; "S"="0.005918093368856362"; "T"="0.005918093368856364" } }
    @{ Row=3; Cells=@{ "E"="3"; "F"="1"; "G"="4.204118999999999"; "H"="12.612357"; "I"="0.01983154129720676"; "J"="0.01983154129720676"; "O"="0.2893586437755394"; "P"="0.2893586437755394"; "Q"="685.297612484711"; "R"="6167.678512362399"; "S"="0.005738427893738349"; "T"="0.005738427893738351" } }
    @{ Row=4; Cells=@{ "E"="3"; "F"="1"; "G"="4.204118999999999"; "H"="12.612357"; "I"="0.01983154129720676"; "J"="0.01983154129720676"; "M"="165.99353"; "N"="497.98059"; "O"="0.294661504941043"; "P"="0.294661504941043"; "Q"="697.8565533500699"; "R"="6280.70898015063"; "S"="0.005843591803935387"; "T"="0.005843591803935387" } }
    @{ Row=5; Cells=@{ "E"="3"; "F"="1"; "G"="4.204118999999999"; "H"="12.612357"; "I"="0.01983154129720676"; "J"="0.01983154129720676"; "M"="66.22673433333334"; "N"="198.680203"; "O"="0.1175616254801657"; "P"="0.1175616254801657"; "Q"="278.425072118719"; "R"="2505.825649068471"; "S"="0.00233142823067666"; "T"="0.00233142823067666" } }
    @{ Row=6; Cells=@{ "I"="0.8539093107807857"; "J"="0.8539093107807858"; "M"="168.1098273333333"; "N"="504.329482"; "O"="0.2984182258032519"; "P"="0.298418225803252"; "Q"="30431.50162618559"; "R"="273883.5146356703"; "S"="0.2548221015200797"; "T"="0.2548221015200798" } }
    @{ Row=7; Cells=@{ "I"="0.8539093107807857"; "J"="0.8539093107807858"; "O"="0.2893586437755394"; "P"="0.2893586437755394"; "S"="0.2470860400748338"; "T"="0.2470860400748338" } }
    @{ Row=8; Cells=@{ "I"="0.8539093107807857"; "J"="0.8539093107807858"; "M"="165.99353"; "N"="497.98059"; "O"="0.294661504941043"; "P"="0.294661504941043"; "Q"="30048.40620123385"; "R"="270435.6558111047"; "S"="0.2516142025978351"; "T"="0.2516142025978351" } }
    @{ Row=9; Cells=@{ "I"="0.8539093107807857"; "J"="0.8539093107807858"; "M"="66.22673433333334"; "N"="198.680203"; "O"="0.1175616254801657"; "P"="0.1175616254801657"; "Q"="11988.46614460938"; "R"="107896.1953014844"; "S"="0.1003869665880371"; "T"="0.1003869665880371" } }
    @{ Row=10; Cells=@{ "G"="26.057747"; "H"="78.173241"; "I"="0.1229188055196976"; "J"="0.1229188055196976"; "M"="168.1098273333333"; "N"="504.329482"; "O"="0.2984182258032519"; "P"="0.298418225803252"; "Q"="4380.563348865685"; "R"="39425.07013979116"; "S"="0.03668121186104313"; "T"="0.03668121186104314" } }
    @{ Row=11; Cells=@{ "G"="26.057747"; "H"="78.173241"; "I"="0.1229188055196976"; "J"="0.1229188055196976"; "O"="0.2893586437755394"; "P"="0.2893586437755394"; "Q"="4247.575248424377"; "R"="38228.17723581939"; "S"="0.03556761885968899"; "T"="0.03556761885968899" } }
    @{ Row=12; Cells=@{ "G"="26.057747"; "H"="78.173241"; "I"="0.1229188055196976"; "J"="0.1229188055196976"; "M"="165.99353"; "N"="497.98059"; "O"="0.294661504941043"; "P"="0.294661504941043"; "Q"="4325.41740837691"; "R"="38928.75667539219"; "S"="0.03621944021998948"; "T"="0.03621944021998948" } }
    @{ Row=13; Cells=@{ "G"="26.057747"; "H"="78.173241"; "I"="0.1229188055196976"; "J"="0.1229188055196976"; "M"="66.22673433333334"; "N"="198.680203"; "O"="0.1175616254801657"; "P"="0.1175616254801657"; "Q"="1725.719487894214"; "R"="15531.47539104792"; "S"="0.01445053457897601"; "T"="0.01445053457897601" } }
    @{ Row=14; Cells=@{ "G"="0.7081243333333332"; "H"="2.124373"; "I"="0.003340342402309973"; "J"="0.003340342402309974"; "M"="168.1098273333333"; "N"="504.329482"; "O"="0.2984182258032519"; "P"="0.298418225803252"; "Q"="119.0426594071984"; "R"="1071.383934664786"; "S"="0.0009968190532727147"; "T"="0.0009968190532727151" } }
    @{ Row=15; Cells=@{ "G"="0.7081243333333332"; "H"="2.124373"; "I"="0.003340342402309973"; "J"="0.003340342402309974"; "O"="0.2893586437755394"; "P"="0.2893586437755394"; "Q"="115.4286819606345"; "R"="1038.858137645711"; "S"="0.0009665569472783412"; "T"="0.0009665569472783414" } }
    @{ Row=16; Cells=@{ "G"="0.7081243333333332"; "H"="2.124373"; "I"="0.003340342402309973"; "J"="0.003340342402309974"; "M"="165.99353"; "N"="497.98059"; "O"="0.294661504941043"; "P"="0.294661504941043"; "Q"="117.5440577688966"; "R"="1057.89651992007"; "S"="0.0009842703192830355"; "T"="0.0009842703192830358" } }
    @{ Row=17; Cells=@{ "G"="0.7081243333333332"; "H"="2.124373"; "I"="0.003340342402309973"; "J"="0.003340342402309974"; "M"="66.22673433333334"; "N"="198.680203"; "O"="0.1175616254801657"; "P"="0.1175616254801657"; "Q"="46.89676209863544"; "R"="422.070858887719"; "S"="0.0003926960824758819"; "T"="0.000392696082475882" } }
)

foreach ($rowUpdate in $updates) {
    $r = $rowUpdate.Row
    foreach ($col in $rowUpdate.Cells.Keys) {
        $ws.Range("$col$r").Value = [double]$rowUpdate.Cells[$col]
    }
}

